# "Generate Report for Archive"
#
# The localization status report previously showed files as "Ready for
# handoff"; regenerating the report for archive now reflects the files as
# "In Translation" instead. This text shows up in the per-language status
# columns (zh-cn / de-de "Status" column, and the Overview sheet's per-
# language status columns), and because the new text is shorter the
# status columns are re-sized to fit the new content.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (E) / de-de (F) status columns, rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2","F2","E3","F3","E4","F4")) {
    $wsOverview.Range($addr).Value = $newStatus
}

# --- zh-cn sheet: Status column (C), rows 2-4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($addr in @("C2","C3","C4")) {
    $wsZhCn.Range($addr).Value = $newStatus
}

# --- de-de sheet: Status column (C), rows 2-4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($addr in @("C2","C3","C4")) {
    $wsDeDe.Range($addr).Value = $newStatus
}

# --- Re-fit the now-narrower status columns ---
# Shrinking "Ready for handoff" (18 chars) down to "In Translation" (14
# chars) narrows the status columns; resize them to fit the new text.
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
